$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the style from E1 (header style)
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (F2:F18)
$times = @(
    "2021-10-05 10:51:15.877150",
    "2021-10-05 10:51:15.877160",
    "2021-10-05 10:51:15.877163",
    "2021-10-05 10:51:15.877166",
    "2021-10-05 10:51:15.877169",
    "2021-10-05 10:51:15.877172",
    "2021-10-05 10:51:15.877174",
    "2021-10-05 10:51:15.877177",
    "2021-10-05 10:51:15.877179",
    "2021-10-05 10:51:15.877182",
    "2021-10-05 10:51:15.877185",
    "2021-10-05 10:51:15.877187",
    "2021-10-05 10:51:15.877190",
    "2021-10-05 10:51:15.877192",
    "2021-10-05 10:51:15.877194",
    "2021-10-05 10:51:15.877197",
    "2021-10-05 10:51:15.877200"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
